$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 259 (B and D values were recomputed/updated upstream)
$ws.Range("B259").Value = 6243090940000
$ws.Range("D259").Value = 286710429898.6448

# Duplicate formatting (date style etc.) of row 259 down into the three new rows
$ws.Range("A259:D259").Copy()
$ws.Range("A260:D262").PasteSpecial(-4122)

# New row 260 (2023-07-01)
$ws.Range("A260").Value = 45108
$ws.Range("B260").Value = 6355692770000
$ws.Range("C260").Value = 0.04603977376055174
$ws.Range("D260").Value = 292614657222.3744

# New row 261 (2023-08-01)
$ws.Range("A261").Value = 45139
$ws.Range("B261").Value = 6337051350000
$ws.Range("C261").Value = 0.04504991530615922
$ws.Range("D261").Value = 285483626608.2819

# New row 262 (2023-09-01)
$ws.Range("A262").Value = 45170
$ws.Range("B262").Value = 6359425540000
$ws.Range("C262").Value = 0.04329632219391124
$ws.Range("D262").Value = 275339737148.028
